$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the diff. Columns D/E hold numeric-looking text
# (prices / percentages) stored as strings in the source workbook, so we
# force Text number format before assigning -- otherwise Excel would auto-
# convert them into numbers/percentages and mangle the exact formatting
# (trailing zeros, precision, "%" suffix).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '244.27'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.49%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '26.40'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '4.12%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.132'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.22%'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.30%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.464'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.52%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8173'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.02%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8337'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-1.25%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1333'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06997'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '0.56%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02889'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '0.01%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09394'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '0.20%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001529'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.01%'
$ws.Range('B14').Value = 'One'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0005970'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.37%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.006199'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.99%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.648'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '4.24%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.034'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '0.49%'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.183'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '5.82%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.03105'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-2.22%'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-2.26%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.739'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-0.12%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04662'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-1.12%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.10%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001244'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-0.52%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-2.98%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.00009603'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-1.03%'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0001393'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '0.33%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03641'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.54%'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1361'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.36%'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.002621'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.38%'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003442'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-44.50%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008872'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '6.68%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005350'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.94%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.04%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-4.04%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '9.98%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.04%'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.04%'
